$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999998762845288
$ws.Range("A2").Value = 0.99581994833251841
$ws.Range("A3").Value = 0.97711516827561795
$ws.Range("A4").Value = 0.969172387318884
$ws.Range("A5").Value = 0.96170058667752145
$ws.Range("A6").Value = 0.94590149909790511
$ws.Range("A7").Value = 0.94441556135928351
$ws.Range("A8").Value = 0.94061077233455892
$ws.Range("A9").Value = 0.93916066183786673
$ws.Range("A10").Value = 0.93864472260512199
$ws.Range("A11").Value = 0.93849764875719166
$ws.Range("A12").Value = 0.93847466181766825
$ws.Range("A13").Value = 0.94390176315783458
$ws.Range("A14").Value = 0.94630137421641836
$ws.Range("A15").Value = 0.94558422269325271
$ws.Range("A16").Value = 0.9453772221725083
$ws.Range("A17").Value = 0.94604359682801453
$ws.Range("A18").Value = 0.94720570098301415
$ws.Range("A19").Value = 0.99476682347902057
$ws.Range("A20").Value = 0.98765004563335013
$ws.Range("A21").Value = 0.98625160119640287
$ws.Range("A22").Value = 0.98498709936966411
$ws.Range("A23").Value = 0.97007997410915969
$ws.Range("A24").Value = 0.95705866710879539
$ws.Range("A25").Value = 0.95060161732037685
$ws.Range("A26").Value = 0.94401482869485387
$ws.Range("A27").Value = 0.94016667594426784
$ws.Range("A28").Value = 0.92561130044084206
$ws.Range("A29").Value = 0.91544837790558597
$ws.Range("A30").Value = 0.91057065387551128
$ws.Range("A31").Value = 0.90998246058034327
$ws.Range("A32").Value = 0.90830323370302934
$ws.Range("A33").Value = 0.90778323998870303
